$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "41.187.28"
$ws.Cells.Item(3,4).Value = "2.141.05"
$ws.Cells.Item(3,5).Value = "  -3.38%  "
$ws.Cells.Item(4,5).Value = "  -0.02%  "
$ws.Cells.Item(5,4).Value = "'234.83"
$ws.Cells.Item(5,5).Value = "  -2.83%  "
$ws.Cells.Item(6,5).Value = "  -4.59%  "
$ws.Cells.Item(7,4).Value = "'68.89"
$ws.Cells.Item(7,5).Value = "  -5.70%  "
$ws.Cells.Item(8,5).Value = "  +0.03%  "
$ws.Cells.Item(9,5).Value = "  -7.01%  "
$ws.Cells.Item(10,4).Value = "'38.21"
$ws.Cells.Item(10,5).Value = "  -10.24%  "
$ws.Cells.Item(11,4).Value = "'0.0892"
$ws.Cells.Item(11,5).Value = "  -7.09%  "
$ws.Cells.Item(12,4).Value = "'53.47"
$ws.Cells.Item(12,5).Value = "  -6.97%  "
$ws.Cells.Item(13,4).Value = "'0.0994"
$ws.Cells.Item(13,5).Value = "  -4.12%  "
$ws.Cells.Item(14,5).Value = "  -6.63%  "
$ws.Cells.Item(15,4).Value = "2.460.33"
$ws.Cells.Item(15,5).Value = "  -3.43%  "
$ws.Cells.Item(16,4).Value = "'14.32"
$ws.Cells.Item(16,5).Value = "  +0.25%  "
$ws.Cells.Item(17,4).Value = "2.134.14"
$ws.Cells.Item(17,5).Value = "  -3.09%  "
$ws.Cells.Item(18,5).Value = "  -7.56%  "
$ws.Cells.Item(19,4).Value = "41.050.61"
$ws.Cells.Item(19,5).Value = "  -2.07%  "
$ws.Cells.Item(20,4).Value = "0.0₃0992"
$ws.Cells.Item(20,5).Value = "  -7.79%  "
$ws.Cells.Item(21,4).Value = "'69.14"
$ws.Cells.Item(21,5).Value = "  -4.93%  "
$ws.Cells.Item(22,5).Value = "  -7.88%  "
$ws.Cells.Item(23,2).Value = "BitcoinCash"
$ws.Cells.Item(23,3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(23,4).Value = "'223.53"
$ws.Cells.Item(23,5).Value = "  -3.06%  "
$ws.Cells.Item(24,2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(24,3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(24,4).Value = "'9.46"
$ws.Cells.Item(24,5).Value = "  -12.25%  "
$ws.Cells.Item(25,5).Value = "  -0.02%  "
$ws.Cells.Item(26,4).Value = "'1.88"
$ws.Cells.Item(26,5).Value = "  -9.74%  "
$ws.Cells.Item(27,5).Value = "  -11.19%  "
$ws.Cells.Item(28,4).Value = "'3.37"
$ws.Cells.Item(28,5).Value = "  -8.29%  "
$ws.Cells.Item(29,5).Value = "  -6.05%  "
$ws.Cells.Item(30,5).Value = "  -2.46%  "
$ws.Cells.Item(31,4).Value = "'168.41"
$ws.Cells.Item(31,5).Value = "  +0.14%  "
$ws.Cells.Item(32,4).Value = "'19.50"
$ws.Cells.Item(32,5).Value = "  -4.87%  "
$ws.Cells.Item(33,4).Value = "'30.84"
$ws.Cells.Item(33,5).Value = "  +2.94%  "
$ws.Cells.Item(34,4).Value = "'0.0749"
$ws.Cells.Item(34,5).Value = "  -5.82%  "
$ws.Cells.Item(35,5).Value = "  -11.66%  "
$ws.Cells.Item(36,5).Value = "  -5.17%  "
$ws.Cells.Item(37,4).Value = "'0.0999"
$ws.Cells.Item(37,5).Value = "  -9.07%  "
$ws.Cells.Item(38,5).Value = "  -5.08%  "
$ws.Cells.Item(39,4).Value = "'0.0279"
$ws.Cells.Item(39,5).Value = "  -7.65%  "
$ws.Cells.Item(40,5).Value = "  -4.63%  "
$ws.Cells.Item(41,4).Value = "'11.57"
$ws.Cells.Item(41,5).Value = "  -16.95%  "
$ws.Cells.Item(42,5).Value = "  -7.16%  "
$ws.Cells.Item(43,4).Value = "'56.94"
$ws.Cells.Item(43,5).Value = "  -13.10%  "
$ws.Cells.Item(44,4).Value = "'0.185"
$ws.Cells.Item(44,5).Value = "  -6.79%  "
$ws.Cells.Item(45,4).Value = "'8.14"
$ws.Cells.Item(45,5).Value = "  -7.71%  "
$ws.Cells.Item(46,5).Value = "  -6.19%  "
$ws.Cells.Item(47,4).Value = "'96.40"
$ws.Cells.Item(47,5).Value = "  -8.43%  "
$ws.Cells.Item(48,5).Value = "  -4.77%  "
$ws.Cells.Item(49,5).Value = "  -5.92%  "
$ws.Cells.Item(50,5).Value = "  -12.00%  "
$ws.Cells.Item(51,5).Value = "  -3.57%  "
